$wb = $excel.ActiveWorkbook

# Loan RBI, Variable Instalments:
# Insert a new (blank) column at N on the "Repayment Schedule" sheet, which
# shifts the existing "Late" / "Heading" / "Outstanding" columns one slot to
# the right (N->O, O->P, P->Q).
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet/tab, with cell R7 selected.
$wsSchedule.Activate()
$wsSchedule.Range("R7").Select()
